# Scheduled-runner refresh of market-price-derived columns (H:N) across
# several Leve-profit sheets. Values come from an external price feed and
# are written as plain numbers (no formulas in this workbook). Some rows
# only ever populate one of the NQ/HQ profit columns (M vs N) depending on
# which is profitable, so a handful of cells are cleared out entirely and
# a few others appear for the first time - handled below via $null resets
# and plain value assignments.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Range("H18").Value = 1159
$ws.Range("I18").Value = 1188.8889
$ws.Range("J18").Value = 890
$ws.Range("K18").Value = 1188.8889
$ws.Range("L18").Value = 890
$ws.Range("M18").Value = -904.8888999999999
$ws.Range("N18").Value = -1458

# Row 40
$ws.Range("H40").Value = 7011.7646
$ws.Range("I40").Value = 7011.7646
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 7011.7646
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = $null
$ws.Range("N40").Value = -6836.7646

# Row 129
$ws.Range("H129").Value = 740.8182
$ws.Range("I129").Value = 317.625
$ws.Range("J129").Value = 982.6429000000001
$ws.Range("K129").Value = 952.875
$ws.Range("L129").Value = 2947.9287
$ws.Range("M129").Value = 4047.125
$ws.Range("N129").Value = -12947.9287

# Row 138
$ws.Range("H138").Value = 6285.3335
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 6285.3335
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = $null
$ws.Range("M138").Value = 18856.0005
$ws.Range("N138").Value = -29136.0005

$ws = $wb.Worksheets.Item("ARM")
# Row 6
$ws.Range("H6").Value = 7000
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 7000
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 7000
$ws.Range("M6").Value = $null
$ws.Range("N6").Value = -7346

# Row 61
$ws.Range("H61").Value = 785.6429000000001
$ws.Range("I61").Value = 761.6579
$ws.Range("J61").Value = 1013.5
$ws.Range("K61").Value = 761.6579
$ws.Range("L61").Value = 1013.5
$ws.Range("M61").Value = -549.6579
$ws.Range("N61").Value = -1437.5

# Row 132
$ws.Range("H132").Value = 2361.3142
$ws.Range("I132").Value = 2514.16
$ws.Range("J132").Value = 1979.2
$ws.Range("K132").Value = 7542.48
$ws.Range("L132").Value = 5937.6
$ws.Range("M132").Value = -5012.48
$ws.Range("N132").Value = -10997.6

# Row 136
$ws.Range("H136").Value = 785.6429000000001
$ws.Range("I136").Value = 761.6579
$ws.Range("J136").Value = 1013.5
$ws.Range("K136").Value = 2284.9737
$ws.Range("L136").Value = 3040.5
$ws.Range("M136").Value = 265.0263
$ws.Range("N136").Value = -8140.5

$ws = $wb.Worksheets.Item("BSM")
# Row 54
$ws.Range("H54").Value = 13765.5
$ws.Range("I54").Value = 5555.643
$ws.Range("J54").Value = 42500
$ws.Range("K54").Value = 5555.643
$ws.Range("L54").Value = 42500
$ws.Range("M54").Value = -5071.643
$ws.Range("N54").Value = -43468

# Row 137
$ws.Range("H137").Value = 43320
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 43320
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 43320
$ws.Range("N137").Value = -53520

$ws = $wb.Worksheets.Item("CRP")
# Row 51
$ws.Range("H51").Value = 17921.334
$ws.Range("I51").Value = 500
$ws.Range("J51").Value = 20099
$ws.Range("K51").Value = 500
$ws.Range("L51").Value = 20099
$ws.Range("M51").Value = 236
$ws.Range("N51").Value = -21571

# Row 53
$ws.Range("H53").Value = 20625
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 20625
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = $null
$ws.Range("M53").Value = 20625
$ws.Range("N53").Value = -21839

# Row 61
$ws.Range("H61").Value = 17921.334
$ws.Range("I61").Value = 500
$ws.Range("J61").Value = 20099
$ws.Range("K61").Value = 500
$ws.Range("L61").Value = 20099
$ws.Range("M61").Value = -152
$ws.Range("N61").Value = -20795

# Row 118
$ws.Range("H118").Value = 40742
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 40742
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 40742
$ws.Range("N118").Value = -44056

$ws = $wb.Worksheets.Item("CUL")
# Row 18
$ws.Range("H18").Value = 724.3929000000001
$ws.Range("I18").Value = 659.9091
$ws.Range("J18").Value = 766.1177
$ws.Range("K18").Value = 1979.7273
$ws.Range("L18").Value = 2298.3531
$ws.Range("M18").Value = -1810.7273
$ws.Range("N18").Value = -2636.3531

# Row 38
$ws.Range("H38").Value = 371.54544
$ws.Range("I38").Value = 668.3333
$ws.Range("J38").Value = 260.25
$ws.Range("K38").Value = 2004.9999
$ws.Range("L38").Value = 780.75
$ws.Range("M38").Value = -1657.9999
$ws.Range("N38").Value = -1474.75

# Row 80
$ws.Range("H80").Value = 2773.1462
$ws.Range("I80").Value = 2670
$ws.Range("J80").Value = 2794.3823
$ws.Range("K80").Value = 8010
$ws.Range("L80").Value = 8383.1469
$ws.Range("M80").Value = -7074
$ws.Range("N80").Value = -10255.1469

# Row 83
$ws.Range("H83").Value = 2773.1462
$ws.Range("I83").Value = 2670
$ws.Range("J83").Value = 2794.3823
$ws.Range("K83").Value = 24030
$ws.Range("L83").Value = 25149.4407
$ws.Range("M83").Value = -19350
$ws.Range("N83").Value = -34509.44070000001

# Row 86
$ws.Range("H86").Value = 568.3333
$ws.Range("I86").Value = 202
$ws.Range("J86").Value = 751.5
$ws.Range("K86").Value = 606
$ws.Range("L86").Value = 2254.5
$ws.Range("M86").Value = 580
$ws.Range("N86").Value = -4626.5

# Row 89
$ws.Range("H89").Value = 568.3333
$ws.Range("I89").Value = 202
$ws.Range("J89").Value = 751.5
$ws.Range("K89").Value = 1818
$ws.Range("L89").Value = 6763.5
$ws.Range("M89").Value = 4110
$ws.Range("N89").Value = -18619.5

# Row 92
$ws.Range("H92").Value = 827.7143
$ws.Range("I92").Value = 572.75
$ws.Range("J92").Value = 1167.6666
$ws.Range("K92").Value = 1718.25
$ws.Range("L92").Value = 3502.9998
$ws.Range("M92").Value = -470.25
$ws.Range("N92").Value = -5998.9998

# Row 107
$ws.Range("H107").Value = 17175.5
$ws.Range("I107").Value = 25437.25
$ws.Range("J107").Value = 652
$ws.Range("K107").Value = 76311.75
$ws.Range("L107").Value = 1956
$ws.Range("M107").Value = -74391.75
$ws.Range("N107").Value = -5796

# Row 113
$ws.Range("H113").Value = 1163.1904
$ws.Range("I113").Value = 1361.0625
$ws.Range("J113").Value = 530
$ws.Range("K113").Value = 4083.1875
$ws.Range("L113").Value = 1590
$ws.Range("M113").Value = -1913.1875
$ws.Range("N113").Value = -5930

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 1951.76
$ws.Range("I132").Value = 2046.579
$ws.Range("J132").Value = 1651.5
$ws.Range("K132").Value = 6139.737
$ws.Range("L132").Value = 4954.5
$ws.Range("M132").Value = -3609.737
$ws.Range("N132").Value = -10014.5

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1030.5
$ws.Range("I22").Value = 1030.5
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1030.5
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -735.5

# Row 27
$ws.Range("H27").Value = 1030.5
$ws.Range("I27").Value = 1030.5
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 1030.5
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -923.5

# Row 46
$ws.Range("H46").Value = 1878.2858
$ws.Range("I46").Value = 2974.5
$ws.Range("J46").Value = 416.66666
$ws.Range("K46").Value = 2974.5
$ws.Range("L46").Value = 416.66666
$ws.Range("M46").Value = -2786.5
$ws.Range("N46").Value = -792.66666

# Row 100
$ws.Range("H100").Value = 2992.6072
$ws.Range("I100").Value = 2896.5
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 2896.5
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -2355.5
$ws.Range("N100").Value = -4082
